# Generate Report for Handback
#
# - Overview sheet: status columns (zh-cn / de-de) move from "Ready for
#   handoff" to "Handed back: in sync with en-US"; those columns (and the
#   matching per-language "Latest Target File" columns) are widened to fit
#   the longer text.
# - zh-cn / de-de sheets: the "Latest Target File" (I) and "Latest Handback
#   File" (J) cells are now populated (they link to the source .md doc and
#   the generated handback .xlf respectively) and the "Latest Handback
#   DateTime" (K) timestamps are stamped with real values.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$wideWidth = 29.9777047293527 - 0.8333333333333334   # -> stored col width ~30
$fullWidth = 40 - 0.8333333333333334                  # -> stored col width 40
$hyperlinkBlue = 15570276                              # RGB(100,149,237)

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Columns.Item(3).ColumnWidth = $wideWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $fullWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $fullWidth

$i2 = $wsZhCn.Range("I2")
$i2.Value = "50e99dc3-41f3-47c3-b97a-a4bcc6e40e4e.md"
$i2.Font.Underline = $true
$i2.Font.Color = $hyperlinkBlue
$wsZhCn.Hyperlinks.Add($i2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/15303f26ca62f5df54879773b36363fb9d8e2300/e2e/50e99dc3-41f3-47c3-b97a-a4bcc6e40e4e.md", "", "", "50e99dc3-41f3-47c3-b97a-a4bcc6e40e4e.md")
$wsZhCn.Range("J2").Value = "50e99dc3-41f3-47c3-b97a-a4bcc6e40e4e.7afd2651fcb03608e1b82ffcd3ded2ba4660e266.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-30 22:52:17"

$i3 = $wsZhCn.Range("I3")
$i3.Value = "bc5a4f33-aa25-4aa1-92cf-2c9dbcbe09ec.md"
$i3.Font.Underline = $true
$i3.Font.Color = $hyperlinkBlue
$wsZhCn.Hyperlinks.Add($i3, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/15303f26ca62f5df54879773b36363fb9d8e2300/e2e/bc5a4f33-aa25-4aa1-92cf-2c9dbcbe09ec.md", "", "", "bc5a4f33-aa25-4aa1-92cf-2c9dbcbe09ec.md")
$wsZhCn.Range("J3").Value = "bc5a4f33-aa25-4aa1-92cf-2c9dbcbe09ec.dbcf3e5de7bb76c7d88afd670220425df25e0728.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-30 22:52:17"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Columns.Item(3).ColumnWidth = $wideWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $fullWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $fullWidth

$i2de = $wsDeDe.Range("I2")
$i2de.Value = "50e99dc3-41f3-47c3-b97a-a4bcc6e40e4e.md"
$i2de.Font.Underline = $true
$i2de.Font.Color = $hyperlinkBlue
$wsDeDe.Hyperlinks.Add($i2de, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/15303f26ca62f5df54879773b36363fb9d8e2300/e2e/50e99dc3-41f3-47c3-b97a-a4bcc6e40e4e.md", "", "", "50e99dc3-41f3-47c3-b97a-a4bcc6e40e4e.md")
$wsDeDe.Range("J2").Value = "50e99dc3-41f3-47c3-b97a-a4bcc6e40e4e.7afd2651fcb03608e1b82ffcd3ded2ba4660e266.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-30 22:52:24"

$i3de = $wsDeDe.Range("I3")
$i3de.Value = "bc5a4f33-aa25-4aa1-92cf-2c9dbcbe09ec.md"
$i3de.Font.Underline = $true
$i3de.Font.Color = $hyperlinkBlue
$wsDeDe.Hyperlinks.Add($i3de, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/15303f26ca62f5df54879773b36363fb9d8e2300/e2e/bc5a4f33-aa25-4aa1-92cf-2c9dbcbe09ec.md", "", "", "bc5a4f33-aa25-4aa1-92cf-2c9dbcbe09ec.md")
$wsDeDe.Range("J3").Value = "bc5a4f33-aa25-4aa1-92cf-2c9dbcbe09ec.dbcf3e5de7bb76c7d88afd670220425df25e0728.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-30 22:52:24"
